$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - ARE_GUI_1 (unchanged content, kept for reference)
# A7/B7/D7/E7 unchanged; C7 changes value (text) from index 45 -> "ARE start file: start.bat (start.sh - Linux)"
$ws.Range("C7").Value = "ARE start file: start.bat (start.sh - Linux)"

# Row 8 - was ARE_GUI_2 "Deploy"/stop-test details, becomes ARE_GUI_2 "Stop model"
$ws.Range("A8").Value = "ARE_GUI_2"
$ws.Range("B8").Value = "Stop model"
$ws.Range("C8").Value = "Test ARE_GUI_1"
$ws.Range("D8").Value = "1. Execute test ARE_GUI_1`n2. Right click on ARE GUI background panel`n3. Click on 'Stop Model' button"
$ws.Range("E8").Value = "1. The model must stop running`n2. CameraMouse: The LED must be off and the ARE GUI panel empty, the video frame window must be closed"

# Row 9 - becomes ARE_GUI_3 "Start model"
$ws.Range("A9").Value = "ARE_GUI_3"
$ws.Range("B9").Value = "Start model"
$ws.Range("C9").Value = "Test ARE_GUI_2"
$ws.Range("D9").Value = "1. Execute test ARE_GUI_2`n2. Right click on ARE GUI background panel`n3. Click on 'Start Model' button"
$ws.Range("E9").Value = "1. The model must be started successfully.`n2. CameraMouse: `nThe camera LED must be on and the camera frames be visible in a dedicated video frame window.`n(The mouse cursor must move according to head movements, if mouse control is enabled)"

# Row 10 - becomes ARE_GUI_4 "Pause model/Start model"
$ws.Range("A10").Value = "ARE_GUI_4"
$ws.Range("B10").Value = "Pause model/Start model"
$ws.Range("C10").Value = "Test ARE_GUI_1"
$ws.Range("D10").Value = "1. Execute test ARE_GUI_1`n2. Right click on ARE GUI background panel`n3. Click on 'Pause Model' button`n4. Click on 'Start Model' button"
$ws.Range("E10").Value = "1. After clicking pause, the model  must pause`nCameraMouse: The LED must be off and the video frame window must be closed`n2. After clicking start, the model must start successfully`nCameraMouse: `nThe camera LED must be on and the camera frames be visible in a dedicated video frame window.`n(The mouse cursor must move according to head movements, if mouse control is enabled)"

# Row 11 - becomes ARE_GUI_5 "Pause model/Stop model"
$ws.Range("A11").Value = "ARE_GUI_5"
$ws.Range("B11").Value = "Pause model/Stop model"
$ws.Range("C11").Value = "Test ARE_GUI_1"
$ws.Range("D11").Value = "1. Execute test ARE_GUI_1`n2. Right click on ARE GUI background panel`n3. Click on 'Pause Model' button`n4. Click on 'Stop Model' button"
$ws.Range("E11").Value = "1. After clicking pause, the model  must pause`nCameraMouse: The LED must be off and the video frame window must be closed`n2. After clicking stop, the model must stop running`nCameraMouse: The LED must be off and the ARE GUI panel empty, the video frame window must be closed"

# Row 12 - becomes ARE_GUI_6 "Test start of all demo models on windows"
$ws.Range("A12").Value = "ARE_GUI_6"
$ws.Range("B12").Value = "Test start of all demo models on windows"
$ws.Range("C12").Value = "Windows OS`nARE start file: start.bat (start.sh - Linux)`ndemomenu.acs (Should be default autostart model)"
$ws.Range("D12").Value = "1. Execute ARE start file`n2. Select each menu entry of the demo menu sequentially (including submenus)"
$ws.Range("E12").Value = "Each demo model should start successfully "

# Row 13 - becomes ARE_GUI_7 "Stresstest Start model"
$ws.Range("A13").Value = "ARE_GUI_7"
$ws.Range("B13").Value = "Stresstest Start model"
$ws.Range("C13").Value = "Test ARE_GUI_2"
$ws.Range("D13").Value = "1. Execute Test ARE_GUI_2 by clicking 10 times onto 'Start' button as fast as possible"
$ws.Range("E13").Value = "The model must be started 10 times sequentially and successfully`nThe last model start must have a clean state and must not have orphaned GUI elements in the ARE GUI  panel`nThe ARE must not crash"

# Update the selected cell in the sheet view
$ws.Range("C9").Select()
